$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8: change the time-range label, add the trigger-minutes JSON into column C,
# and switch the message type from "daily" to "hourly" (matching the pattern used in row 25).
$ws.Range("B8").Value = "10:55 - 10:59"
$ws.Range("C8").Value = '["01", "04", "08", "13", "17", "22", "26", "30"]'
$ws.Range("D8").Value = "часовой посыл"

# Row 9: same treatment - update the paired time-range label, add the trigger JSON,
# and switch the type to "hourly" as well.
$ws.Range("B9").Value = "11:00 - 11:04"
$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = '["01", "04", "08", "13", "17", "22", "26", "30"]'
$ws.Range("D9").Value = "часовой посыл"

# Update the visible selection/scroll position to match the saved view.
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("B12").Select()
